# Update results sheets (2025, 2030, 2035, 2040, 2045, 2050) with new
# server-computed headers and values as described in commit "ADD results
# from server".
#
# New header layout (columns A..O):
#   eb, gb, hp, st, wi, ieh, chp, ac, ab_ct, ab_hp, cp_ct, cp_hp, ttes, btes, ites
#
# Row 2 holds the corresponding numeric results for each sheet/year.

$wb = $excel.ActiveWorkbook

$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")

$data = @{
    1 = @(3906.399109145206, 0, 48353.76274462014, 0, 289724.0114301849, 9433.134471502228, 0, 2534.277928792104, 0, 0, 0, 0, 0, 2367.37219622158, 1995.762462679798)
    2 = @(6991.052031681918, 0, 197913.7502057619, 0, 289724.0114301849, 16452.51445364119, 0, 8194.52068131253, 0, 0, 0, 0, 0, 7543.193583625169, 6257.586732772244)
    3 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 12888.04225687751, 9263.466444480218)
    4 = @(31236.29455387744, 0, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 14045.89200932069, 9263.466444480218)
    5 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16880.86083028515, 10096.08766803326)
    6 = @(38906.8534480406, 193.0947398408091, 292247.2772138842, 0, 289724.0114301849, 16595.10705160327, 0, 12131.91920790125, 0, 0, 0, 0, 0, 16880.86083028515, 10096.08766803326)
}

for ($sheetIndex = 1; $sheetIndex -le 6; $sheetIndex++) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $values = $data[$sheetIndex]

    for ($col = 1; $col -le 15; $col++) {
        $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
        $ws.Cells.Item(2, $col).Value = $values[$col - 1]
    }
}
